$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") from row 2 to row 34: 45625 -> 45626
for ($r = 2; $r -le 34; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45625) {
        $cell.Value = 45626
    }
}
